$d = $word.ActiveDocument

# --- Change 2 (done first to avoid name clash): remove the old _GoBack bookmark ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Change 1: insert three new paragraphs right after the "Summary" heading ---
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("Summary", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$summaryPara = $findRange.Paragraphs(1)
$insertPos = $summaryPara.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

$newParasXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
'<w:p><w:pPr><w:pStyle w:val="Footer"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:b/><w:bCs/><w:smallCaps/><w:spacing w:val="60"/><w:sz w:val="28"/></w:rPr></w:pPr></w:p>' + `
'<w:p><w:pPr><w:pStyle w:val="Footer"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">I am a senior software developer of nearly 17 years who gets up every day ready to solve more problems with technology. </w:t></w:r></w:p>' + `
'<w:p><w:pPr><w:pStyle w:val="Footer"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($newParasXml)

function Replace-ParagraphByAnchor([string]$anchorText, [string]$newParaInnerXml, [string]$paraAttrs) {
    $fr = $d.Content
    $fr.Find.ClearFormatting()
    $ok = $fr.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Anchor text not found: $anchorText"
    }
    $para = $fr.Paragraphs(1)
    $start = $para.Range.Start
    $end = $para.Range.End - 1
    $r = $d.Range($start, $end)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p' + $paraAttrs + '>' + $newParaInnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# --- Change 3: split "...collected in smooth, intuitive way" so that
#     "intuitive way" starts a new run carrying <w:lastRenderedPageBreak/> ---
$hunk3Attrs = ' w14:paraId="6E90AE9E" w14:textId="3764151B" w:rsidR="00B40F54" w:rsidRDefault="00B40F54" w:rsidP="00C47CB5"'
$hunk3Inner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
'<w:r w:rsidRPr="00670070"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/></w:rPr><w:t>Business value:</w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> Data services support the primary product offering of my client to allow their clients to reimburse their drivers using accurate mileage data collected in smooth, </w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>intuitive way</w:t></w:r>'
Replace-ParagraphByAnchor "intuitive way" $hunk3Inner $hunk3Attrs

# --- Change 4: remove <w:lastRenderedPageBreak/> before "Mentored new developers" ---
$hunk4Attrs = ' w14:paraId="02C41766" w14:textId="4C2A067E" w:rsidR="00C47CB5" w:rsidRPr="00ED3812" w:rsidRDefault="00C47CB5" w:rsidP="008252A5"'
$hunk4Inner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
'<w:r w:rsidRPr="00ED3812"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Mentored new developers</w:t></w:r>' + `
'<w:r w:rsidR="007D3A98"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' + `
'<w:r w:rsidR="00ED3812" w:rsidRPr="00ED3812"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>interns</w:t></w:r>' + `
'<w:r w:rsidR="007D3A98"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>,</w:t></w:r>' + `
'<w:r w:rsidRPr="00ED3812"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> and testers through frequent collaboration</w:t></w:r>'
Replace-ParagraphByAnchor "Mentored new developers" $hunk4Inner $hunk4Attrs

# --- Change 5: split "... with Visual Studio 2005 and SQL Server 2000/2005." so the
#     second part starts a new run carrying <w:lastRenderedPageBreak/> ---
$hunk5Attrs = ' w14:paraId="673D26BE" w14:textId="77777777" w:rsidR="00B67F06" w:rsidRPr="001F5DA4" w:rsidRDefault="00B67F06" w:rsidP="00B67F06"'
$hunk5Inner = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
'<w:r w:rsidRPr="001F5DA4"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Gained experience in .NET 2.0 Windows and </w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>w</w:t></w:r>' + `
'<w:r w:rsidRPr="001F5DA4"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">eb </w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>d</w:t></w:r>' + `
'<w:r w:rsidRPr="001F5DA4"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>evelopment through everyday use of C#</w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> 2.0</w:t></w:r>' + `
'<w:r w:rsidRPr="001F5DA4"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>with Visual Studio 2005 and SQL Server 2000/2005.</w:t></w:r>'
Replace-ParagraphByAnchor "with Visual Studio 2005" $hunk5Inner $hunk5Attrs

# --- Change 6: merge the "Gained limited experience..." / "Photoshop CS..." runs
#     into a single run, dropping <w:lastRenderedPageBreak/> ---
$hunk6Attrs = ' w14:paraId="29BB2C78" w14:textId="77777777" w:rsidR="00B67F06" w:rsidRDefault="00B67F06" w:rsidP="00B67F06"'
$hunk6Inner = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Gained limited experience in Web Services, Atlas, Flash integration (with .NET), Adobe Photoshop CS, and web design</w:t></w:r>'
Replace-ParagraphByAnchor "Gained limited experience" $hunk6Inner $hunk6Attrs

